$wb = $excel.ActiveWorkbook

# --- Remove the "formula" attribute row from the !!_Schema sheet ---
$schema = $wb.Worksheets.Item("!!_Schema")
$schema.Rows.Item(7).Delete()

# --- Remove the "Id (iAF1260 [Ref1])" column from the !!Metabolites sheet ---
$metabolites = $wb.Worksheets.Item("!!Metabolites")
$metabolites.Columns.Item(2).Delete()

# --- Update the ObjTables generation timestamp everywhere it appears ---
$oldDate = "2020-05-29 00:17:39"
$newDate = "2020-06-02 16:10:23"

foreach ($sheetInfo in @(
        @{ Name = "!!_Table of contents"; Cells = @("A1", "A2") },
        @{ Name = "!!_Schema"; Cells = @("A1") },
        @{ Name = "!!Compartments"; Cells = @("A1") },
        @{ Name = "!!Metabolites"; Cells = @("A1") },
        @{ Name = "!!Reactions"; Cells = @("A1") },
        @{ Name = "!!References"; Cells = @("A1") },
        @{ Name = "!!Regulations"; Cells = @("A1") }
    )) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)
    foreach ($addr in $sheetInfo.Cells) {
        $cell = $ws.Range($addr)
        $cell.Value = $cell.Value.Replace($oldDate, $newDate)
    }
}
